$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Give each other enough time*") {
        $p.Range.HighlightColorIndex = 7
    }
}
